$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.064.55"
$ws.Range("E2").Value = "  -1.97%  "
$ws.Range("D3").Value = "2.103.42"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  -0.70%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "344.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("E7").Value = "  -1.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4419"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09452"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.177"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.12%  "
$ws.Range("D13").Value = "2.109.71"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.736"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.084"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("E16").Value = "  +1.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001168"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.13%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.186"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.53%  "
$ws.Range("D23").Value = "30.156.88"
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.334"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.80%  "
$ws.Range("D26").Value = "2.355.95"
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("E32").Value = "  -1.91%  "
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.251"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.948"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.171"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02576"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06768"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2279"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6950"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.17%  "
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.300"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6671"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.282"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.638"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000352"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.223"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07188"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.85%  "
